# ZBP_06_home_office.xlsx update: add new weekly survey columns
# (week of 13.-19. 9. 2021 and 20.-26. 9. 2021) to both sheets, and
# bump the "aktualizace" (last updated) date in each sheet's footer title.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "data": percentages, new columns BL (13.-19. 9. 2021) and
# BM (20.-26. 9. 2021).
# ---------------------------------------------------------------
$wsData = $wb.Worksheets.Item("data")

# Header row 1 - copy the style of the previous header cell (BK1) so the
# new header cells get the same bold/centered/bordered look.
$wsData.Range("BK1").Copy()
$wsData.Range("BL1:BM1").PasteSpecial(-4122)
$wsData.Range("BL1").Value = "13.–19. 9. 2021"
$wsData.Range("BM1").Value = "20.–26. 9. 2021"

# Data rows 2-77: [row, BL value, BM value]
$dataRows = @(
    @(2, 0.78, 0.8),
    @(3, 0.07000000000000001, 0.06),
    @(4, 0.07000000000000001, 0.07000000000000001),
    @(5, 0.08, 0.07000000000000001),
    @(6, 0.76, 0.79),
    @(7, 0.05, 0.05),
    @(8, 0.06, 0.06),
    @(9, 0.13, 0.1),
    @(10, 0.48, 0.47),
    @(11, 0.21, 0.17),
    @(12, 0.24, 0.3),
    @(13, 0.07000000000000001, 0.06),
    @(14, 0.84, 0.88),
    @(15, 0.04, 0.05),
    @(16, 0.04, 0.02),
    @(17, 0.08, 0.05),
    @(18, 0.86, 0.85),
    @(19, 0.015, 0.02),
    @(20, 0.01, 0.03),
    @(21, 0.115, 0.1),
    @(22, 0.84, 0.84),
    @(23, 0.09, 0.1),
    @(24, 0.04, 0.04),
    @(25, 0.03, 0.02),
    @(26, 0.82, 0.84),
    @(27, 0.05, 0.05),
    @(28, 0.05, 0.05),
    @(29, 0.08, 0.06),
    @(30, 0.82, 0.85),
    @(31, 0.06, 0.05),
    @(32, 0.015, 0.03),
    @(33, 0.105, 0.07000000000000001),
    @(34, 0.71, 0.71),
    @(35, 0.07000000000000001, 0.06),
    @(36, 0.15, 0.17),
    @(37, 0.07000000000000001, 0.06),
    @(38, 0.64, 0.65),
    @(39, 0.12, 0.14),
    @(40, 0.15, 0.14),
    @(41, 0.09, 0.07000000000000001),
    @(42, 0.74, 0.75),
    @(43, 0.1, 0.09),
    @(44, 0.08, 0.11),
    @(45, 0.08, 0.05),
    @(46, 0.83, 0.85),
    @(47, 0.05, 0.05),
    @(48, 0.05, 0.05),
    @(49, 0.07000000000000001, 0.05),
    @(50, 0.67, 0.6899999999999999),
    @(51, 0.08, 0.08),
    @(52, 0.09, 0.08),
    @(53, 0.16, 0.15),
    @(54, 0.79, 0.8100000000000001),
    @(55, 0.04, 0.04),
    @(56, 0.03, 0.04),
    @(57, 0.14, 0.11),
    @(58, 0.73, 0.74),
    @(59, 0.1, 0.09),
    @(60, 0.13, 0.13),
    @(61, 0.04, 0.04),
    @(62, 0.79, 0.8100000000000001),
    @(63, 0.07000000000000001, 0.06),
    @(64, 0.06, 0.07000000000000001),
    @(65, 0.08, 0.06),
    @(66, 0.84, 0.86),
    @(67, 0.015, 0.015),
    @(68, 0.02, 0.02),
    @(69, 0.125, 0.105),
    @(70, 0.76, 0.77),
    @(71, 0.08, 0.08),
    @(72, 0.09, 0.09),
    @(73, 0.07000000000000001, 0.06),
    @(74, 0.71, 0.74),
    @(75, 0.13, 0.11),
    @(76, 0.09, 0.11),
    @(77, 0.07000000000000001, 0.04)
)

foreach ($r in $dataRows) {
    $rowNum = $r[0]
    $wsData.Cells.Item($rowNum, 64).Value = $r[1]
    $wsData.Cells.Item($rowNum, 65).Value = $r[2]
}

# Footer title row 78 - bump the update date.
$wsData.Range("A78").Value = "Život během pandemie, Home office, % respondentů celkově a ve skupinách, aktualizace 6. 10. 2021"

# ---------------------------------------------------------------
# Sheet "pocetR": sample sizes, new columns BK (13.-19. 9. 2021) and
# BL (20.-26. 9. 2021).
# ---------------------------------------------------------------
$wsPocet = $wb.Worksheets.Item("pocetR")

$wsPocet.Range("BJ1").Copy()
$wsPocet.Range("BK1:BL1").PasteSpecial(-4122)
$wsPocet.Range("BK1").Value = "13.–19. 9. 2021"
$wsPocet.Range("BL1").Value = "20.–26. 9. 2021"

# Data rows 2-20: [row, BK value, BL value]
$pocetRows = @(
    @(2, 1043, 1043),
    @(3, 290, 290),
    @(4, 100, 100),
    @(5, 269, 269),
    @(6, 144, 144),
    @(7, 92, 92),
    @(8, 503, 503),
    @(9, 263, 263),
    @(10, 130, 130),
    @(11, 147, 147),
    @(12, 265, 265),
    @(13, 593, 593),
    @(14, 185, 185),
    @(15, 207, 207),
    @(16, 153, 153),
    @(17, 683, 683),
    @(18, 396, 396),
    @(19, 403, 403),
    @(20, 244, 244)
)

foreach ($r in $pocetRows) {
    $rowNum = $r[0]
    $wsPocet.Cells.Item($rowNum, 63).Value = $r[1]
    $wsPocet.Cells.Item($rowNum, 64).Value = $r[2]
}

# Footer title row 21 - bump the update date.
$wsPocet.Range("A21").Value = "Život během pandemie, Home office, velikost dotázaného souboru celkově a ve skupinách, aktualizace 6. 10. 2021"

# Row 21 also carries empty placeholder cells across every other column;
# extend that same pattern into the two new columns (format-only paste so
# the cells exist but stay valueless, matching B21:BJ21).
$wsPocet.Range("BI21").Copy()
$wsPocet.Range("BK21:BL21").PasteSpecial(-4122)
